$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "downsampling" / variable display-name annotations in columns A & B ---
# Column B is filled first (top to bottom), establishing the broad category
# ("fejl" = error, "Mister top" = loses top layer), then column A is filled
# with the more specific reasoning per row.

# Column B (rows 14-24)
$ws.Range("B14").Value = "fejl"
$ws.Range("B15").Value = "Mister top"
$ws.Range("B18").Value = "Mister top"
$ws.Range("B19").Value = "fejl"
$ws.Range("B20").Value = "fejl"
$ws.Range("B21").Value = "fejl"
$ws.Range("B22").Value = "fejl"
$ws.Range("B23").Value = "fejl"
$ws.Range("B24").Value = "fejl"

# Column A (row 13, then 15, 14, then 18-24)
$ws.Range("A13").Value = "reason:"
$ws.Range("A15").Value = "fix efter?"
$ws.Range("A14").Value = "Minus"
$ws.Range("A18").Value = "fix efter?"
$ws.Range("A19").Value = "binær"
$ws.Range("A20").Value = "Binær"
$ws.Range("A21").Value = "ens hele vej"
$ws.Range("A22").Value = 0
$ws.Range("A23").Value = 0
$ws.Range("A24").Value = 0

# --- Update the active selection to match the edited workbook ---
$ws.Range("A29").Select()
